$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing totals row (old row 6), pushing it to row 7
$ws.Rows(6).Insert()

# Row 7 (former row 6, the per-week totals): add row label + bold the label and the grand total
$ws.Range("A7").Value = "Total par semaine"
$ws.Range("A7").Font.Bold = $true
$ws.Range("R7").Font.Bold = $true

# Row 1: header row -> bold "Membres" (A1) and new "Total par personne" label (R1, column of row-sums), both bold
$ws.Range("R1").Value = "Total par personne"
$ws.Range("A1").Font.Bold = $true
$ws.Range("R1").Font.Bold = $true

# New row 6: label for the grand-total cell below it, bold
$ws.Range("R6").Value = "Total"
$ws.Range("R6").Font.Bold = $true

# Column widths (closest achievable values to the authored 15.69140625 / 16.69140625)
$ws.Columns("A").ColumnWidth = 14.83
$ws.Columns("R").ColumnWidth = 15.83

# Selection matches the post-edit cursor location recorded in the workbook
[void]$ws.Range("A9").Select()

# Page setup (paper size / orientation) as captured by the diff
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
